# Incluye Enmascaramiento de PRESTAMOANEXO
# - Activa (ACTIVO = "Y") las filas existentes de PRESTAMOANEXO (filas 204-206)
# - Agrega filas nuevas de enmascaramiento para AGVIRTUAL (filas 381-391)
# - Amplia el rango con nombre _FilterDatabase a la nueva extension de datos

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Activar el enmascaramiento para las columnas de PRESTAMOANEXO ---
$ws.Cells.Item(204, 5).Value = "Y"
$ws.Cells.Item(205, 5).Value = "Y"
$ws.Cells.Item(206, 5).Value = "Y"

# --- Nuevas filas: enmascaramiento de AGVIRTUAL ---
$newRows = @(
    @("AGVIRTUAL", "OPER_FRECUENTE",       "NOMBRESOCIO_EMPRESA", "NOMBRE", "Y", "NOMBRE"),
    @("AGVIRTUAL", "PAC_CLIENTE_MAE",      "DES_NOMBRE_CORTO",    "NOMBRE", "Y", "NOMBRE"),
    @("AGVIRTUAL", "PAC_CLIENTE_MAE",      "DES_NOMBRES",         "NOMBRE", "Y", "NOMBRE"),
    @("AGVIRTUAL", "PAC_CLIENTE_MAE",      "DES_APELLIDOS",       "NOMBRE", "Y", "APELLIDO"),
    @("AGVIRTUAL", "PAC_CLIENTE_MAE",      "DES_CELULAR",         "NUMERO", "Y", "TELEFONO"),
    @("AGVIRTUAL", "PAC_CLIENTE_PACINET",  "DES_NOMBRES",         "NOMBRE", "Y", "NOMBRE"),
    @("AGVIRTUAL", "PAC_CLIENTE_PACINET",  "DES_EMAIL",           "CORREO", "Y", "CORREO"),
    @("AGVIRTUAL", "PAC_CLIENTE_PACINET",  "DES_APELLIDOS",       "NOMBRE", "Y", "APELLIDO"),
    @("AGVIRTUAL", "PAC_CLIENTE_PACINET",  "DES_CELULAR",         "NUMERO", "Y", "TELEFONO"),
    @("AGVIRTUAL", "PAC_USUARIO_MAE",      "DES_NOMBRES",         "NOMBRE", "Y", "NOMBRE"),
    @("AGVIRTUAL", "PAC_USUARIO_MAE",      "DES_APELLIDOS",       "NOMBRE", "Y", "APELLIDO")
)

$startRow = 381
$endRow = $startRow + $newRows.Count - 1
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}

# Replica el formato (estilo) usado por el resto de filas de datos (columnas D:F)
# en las filas recien agregadas, copiando desde una celda de referencia ya existente.
$ws.Cells.Item(206, 6).Copy()
$ws.Range("D" + $startRow + ":F" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Ampliar el rango con nombre _FilterDatabase (oculto) hasta la nueva ultima fila ---
$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = '=ENMASCARAR_COLUMNAS!$A$1:$F$391'
